$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40 (shifts existing rows 40-61 down to 41-62),
# inheriting the date style from the row above for column D.
$ws.Rows.Item(40).Insert()

$ws.Range("A40").Value = 3
$ws.Range("B40").Value = "Femacal de La Calera"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44767
$ws.Range("E40").Value = 5
$ws.Range("F40").Value = 100112022
$ws.Range("G40").Value = "Arveja Verde"
$ws.Range("H40").Value = "Perfection"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 45
$ws.Range("K40").Value = 37000
$ws.Range("L40").Value = 38000
$ws.Range("M40").Value = 37556
$ws.Range("N40").Value = "`$/saco 25 kilos"
$ws.Range("O40").Value = "Provincia de Limarí"
$ws.Range("P40").Value = 1502
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"
